$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.848400000000002
$ws.Range("C3").Value = -11.65199999999999
$ws.Range("E8").Value = 16.0407
$ws.Range("E11").Value = 16.226
$ws.Range("A12").Value = -21.6162
$ws.Range("B14").Value = 6.616199999999997
$ws.Range("E14").Value = 17.1022
$ws.Range("E15").Value = 16.413
$ws.Range("E17").Value = 16.40400000000001
$ws.Range("C20").Value = -11.54880000000001
$ws.Range("C25").Value = -13.20550000000001
$ws.Range("B26").Value = 4.056500000000006
$ws.Range("E26").Value = 15.93970000000001
$ws.Range("A27").Value = -21.5368
$ws.Range("C30").Value = -12.28159999999999
$ws.Range("B31").Value = 4.854800000000003
$ws.Range("A32").Value = -21.36330000000001
$ws.Range("B35").Value = 9.146900000000006
$ws.Range("A36").Value = -20.0643
$ws.Range("E36").Value = 16.49820000000002
$ws.Range("B37").Value = 8.509300000000003
$ws.Range("A38").Value = -19.3995
$ws.Range("C44").Value = -12.98879999999999
$ws.Range("B45").Value = 6.5651
$ws.Range("A46").Value = -21.39260000000001
$ws.Range("C47").Value = -12.1855
$ws.Range("B52").Value = 5.181800000000001
$ws.Range("A54").Value = -21.68939999999999
$ws.Range("A55").Value = -22.7794
$ws.Range("A56").Value = -22.1362
$ws.Range("B57").Value = 4.933999999999996
$ws.Range("C58").Value = -13.21810000000001
$ws.Range("E64").Value = 17.41
$ws.Range("A67").Value = -21.57269999999998
$ws.Range("A69").Value = -21.61239999999999
$ws.Range("A72").Value = -22.06620000000001
$ws.Range("C78").Value = -11.15930000000001
$ws.Range("E79").Value = 18.32850000000002
$ws.Range("B81").Value = 6.577599999999999
$ws.Range("A83").Value = -20.50749999999998
$ws.Range("B83").Value = 7.9364
$ws.Range("C84").Value = -13.52329999999999
$ws.Range("A86").Value = -22.0388
$ws.Range("C89").Value = -10.7665
$ws.Range("E89").Value = 17.48640000000001
$ws.Range("A91").Value = -21.57639999999999
$ws.Range("C91").Value = -11.14319999999999
$ws.Range("C92").Value = -11.44049999999999
$ws.Range("A93").Value = -21.2234
$ws.Range("C96").Value = -12.5897
$ws.Range("A99").Value = -20.25119999999998
$ws.Range("B100").Value = 5.516799999999998
$ws.Range("B102").Value = 8.032400000000003
$ws.Range("C102").Value = -13.3321
